{"js": "// The document contains a paragraph built from three separate runs:\n//   <id>  +  p068r_3  +  </id>\n// The edit merges them into a single run whose text is the\n// concatenation \"<id>p068r_3</id>\", carrying the formatting\n// (Courier New / color 7f6000 / sz 18) of the first (\"<id>\") run.\n//\n// Searching for the full logical string finds a Range that spans all\n// three runs; replacing its text with the same concatenated string\n// collapses them into one run using the formatting of the first\n// character of the matched range (i.e. the original \"<id>\" run),\n// which reproduces the diff exactly.\n\nconst results = context.document.body.search(\"<id>p068r_3</id>\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target '<id>p068r_3</id>' text in the document.\");\n}\n\nresults.items[0].insertText(\"<id>p068r_3</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains a paragraph built out of three separate runs:\n#   \"<id>\"  +  \"p068r_3\"  +  \"</id>\"\n# The edit merges them into a single run containing the concatenated\n# text \"<id>p068r_3</id>\", keeping the formatting (Courier New,\n# color 7f6000, sz 18) that the original \"<id>\" run already carried.\n#\n# Approach: delete the text of the 2nd+3rd runs (\"p068r_3</id>\") so\n# only the \"<id>\" run (and its formatting) remains, then insert the\n# same text back right after the (now sole) \"<id>\" run. Word grows\n# that existing run in place, which keeps its original rPr/rsid\n# attributes intact instead of fabricating a brand-new run node (as a\n# plain Find/Replace would).\n\n$d = $word.ActiveDocument\n\n# Step 1: remove \"p068r_3</id>\" (the 2nd and 3rd runs' text).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"p068r_3</id>\"\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find 'p068r_3</id>' in the document.\"\n}\n\n$range = $find.Parent\n$range.Delete()\n\n# Step 2: find the now-isolated \"<id>\" run and append the removed\n# text right after it, so it becomes part of that same run.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"<id>\"\n$find2.MatchCase = $true\n$find2.MatchWildcards = $false\n$found2 = $find2.Execute()\n\nif (-not $found2) {\n    throw \"Could not find '<id>' in the document.\"\n}\n\n$range2 = $find2.Parent\n$range2.Collapse(0)  # wdCollapseEnd\n$range2.InsertAfter(\"p068r_3</id>\")\n"}
